# Se agrego la grafica y como se inicia el proyecto
#
# The contact list is replaced with a new batch of 12 people (dated
# 12/02/2022 - 23/02/2022) inserted at the top of the table, the last
# 14 rows of the old list are dropped, and the first 9 rows of the old
# list survive, shifted down underneath the new entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Drop the trailing rows of the old data set (old rows 11-24) -
#    none of these people remain in the updated sheet.
$ws.Rows("11:24").Delete()

# 2. Make room for the 12 brand-new records at the top of the table
#    (just under the header row), pushing the surviving old rows
#    (old rows 2-10) down to rows 14-22.
$ws.Rows("2:13").Insert()

# Inserting rows copies the formatting of the row above, so clear it
# back to the unstyled look the rest of the data rows use.
$ws.Range("A2:D13").ClearFormats()

# 3. Fill in the new records. Dates are written with a leading
#    apostrophe so they are stored as plain text, matching the rest
#    of the ULTIMO_PAGO column.
$newRecords = @(
    @("Daniela",   "Villamizar",  3218490916, "12/02/2022"),
    @("Juan",      "Botero",      3218748814, "13/02/2022"),
    @("Cristian",  "Solarte",     3148227994, "14/02/2022"),
    @("Julian",    "Aristizabal", 3046145922, "15/02/2022"),
    @("Juan",      "Londoño",     3163610054, "16/02/2022"),
    @("Mauricio",  "Herrera",     3117754781, "17/02/2022"),
    @("Esteban",   "Meneses",     3108017554, "18/02/2022"),
    @("Carlos",    "Paraco",      3188288098, "19/02/2022"),
    @("Alejandra", "Ruiz",        3167357054, "20/02/2022"),
    @("Daniela",   "Bustos",      3164224295, "21/02/2022"),
    @("Laura",     "Hoyos",       3113829197, "22/02/2022"),
    @("Manuela",   "Rojas",       3114244572, "23/02/2022")
)

$r = 2
foreach ($rec in $newRecords) {
    $ws.Cells.Item($r, 1).Value = $rec[0]
    $ws.Cells.Item($r, 2).Value = $rec[1]
    $ws.Cells.Item($r, 3).Value = $rec[2]
    $ws.Cells.Item($r, 4).Value = "'" + $rec[3]
    $r = $r + 1
}
